# "finestra incidenza 7gg centrata su ultimo g"
#
# The 7-day rolling sum (column C, "somma mobile 7gg.") and the related
# incidence-per-100k figure (column D) used to be computed on a window
# CENTERED on the current day (3 days before .. 3 days after). This edit
# moves the window so that it is centered on ("ends on") the last day,
# i.e. it becomes a trailing window (6 days before .. current day).
# Numerically this means every value that used to sit on row r now
# belongs on row r+3: columns C and D simply shift down by 3 rows, the
# first three previously-populated rows become blank, and three more
# rows gain values at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = 1
$colC = 3
$colD = 4
$shift = 3

# Last row actually used on the sheet (based on column A, the date column).
$lastRow = $ws.Cells.Item($ws.Rows.Count, $colA).End(-4162).Row  # xlUp

# Find the first and last rows that currently hold a real (numeric) value
# in column C - this is the populated range that needs to move down.
# Row 1 is the header row and is skipped.
$firstPopulated = -1
$lastPopulated = -1
for ($r = 2; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, $colC).Value2
    if ($v -ne $null -and $v -ne "") {
        if ($firstPopulated -eq -1) { $firstPopulated = $r }
        $lastPopulated = $r
    }
}

if ($firstPopulated -ne -1) {
    # Walk bottom-up so a destination row is never overwritten before it
    # has been read as a source for a later (lower) destination.
    for ($r = $lastPopulated + $shift; $r -ge ($firstPopulated + $shift); $r--) {
        $src = $r - $shift
        $cVal = $ws.Cells.Item($src, $colC).Value2
        $dVal = $ws.Cells.Item($src, $colD).Value2
        $ws.Cells.Item($r, $colC).Value2 = $cVal
        $ws.Cells.Item($r, $colD).Value2 = $dVal
    }

    # The rows that used to start the populated range no longer have a
    # source feeding them, so blank them out.
    for ($r = $firstPopulated; $r -lt ($firstPopulated + $shift); $r++) {
        $ws.Cells.Item($r, $colC).ClearContents()
        $ws.Cells.Item($r, $colD).ClearContents()
    }
}
